# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# Both sheets list the same set of exhibitions, so the same row -> new value
# updates are applied to each (row numbers differ slightly between the two
# sheets because "全部类型" has one extra row near the top of the table).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new F-column value, for the "展览" sheet.
$exhibitionUpdates = @{
    2  = 2561
    5  = 1443
    6  = 1118
    7  = 323
    8  = 528
    13 = 8878
    14 = 384
    15 = 2495
    16 = 252
    21 = 1164
    23 = 2054
    24 = 2123
    26 = 1823
    27 = 240
    30 = 395
    32 = 125
    33 = 198
    37 = 270
    38 = 453
    39 = 1207
    40 = 275
    41 = 45
    43 = 275
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F-column value, for the "全部类型" sheet.
$allTypesUpdates = @{
    2  = 2561
    5  = 1443
    7  = 1118
    8  = 323
    9  = 528
    13 = 8878
    14 = 384
    15 = 2495
    17 = 252
    22 = 1164
    24 = 2054
    25 = 2123
    27 = 1823
    28 = 240
    31 = 395
    33 = 125
    34 = 198
    38 = 270
    39 = 453
    44 = 1208
    46 = 275
    47 = 45
    49 = 275
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
